# Update Name of Algo - refresh imputed values in result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value  = -8.244499999999988
$ws.Range("D6").Value  = -8.296699999999994
$ws.Range("C7").Value  = -11.62359999999999
$ws.Range("A8").Value  = -21.02760000000001
$ws.Range("D9").Value  = -7.899700000000007
$ws.Range("A10").Value = -20.46769999999998
$ws.Range("D10").Value = -6.755699999999994
$ws.Range("A12").Value = -22.77320000000003
$ws.Range("E12").Value = 12.13049999999999
$ws.Range("B13").Value = 5.846699999999998
$ws.Range("A18").Value = -22.33690000000002
$ws.Range("C20").Value = -14.47160000000001
$ws.Range("E20").Value = 12.2039
$ws.Range("E23").Value = 13.8104
$ws.Range("E25").Value = 13.28349999999999
